$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.96715833333333
$ws.Range("H2").Value = 35.901475
$ws.Range("I2").Value = 0.8213389035667636
$ws.Range("J2").Value = 0.8213389035667636
$ws.Range("M2").Value = 57.65261933333333
$ws.Range("N2").Value = 172.957858
$ws.Range("O2").Value = 0.6817060950001529
$ws.Range("P2").Value = 0.6817060950001529
$ws.Range("Q2").Value = 689.9380238933943
$ws.Range("R2").Value = 6209.442215040549
$ws.Range("S2").Value = 0.5599117366222056
$ws.Range("T2").Value = 0.5599117366222056
$ws.Range("G3").Value = 11.96715833333333
$ws.Range("H3").Value = 35.901475
$ws.Range("I3").Value = 0.8213389035667636
$ws.Range("J3").Value = 0.8213389035667636
$ws.Range("O3").Value = 0.1019529789289588
$ws.Range("P3").Value = 0.1019529789289588
$ws.Range("Q3").Value = 103.1841101732778
$ws.Range("R3").Value = 928.6569915594998
$ws.Range("S3").Value = 0.0837379479288764
$ws.Range("T3").Value = 0.0837379479288764
$ws.Range("G4").Value = 11.96715833333333
$ws.Range("H4").Value = 35.901475
$ws.Range("I4").Value = 0.8213389035667636
$ws.Range("J4").Value = 0.8213389035667636
$ws.Range("M4").Value = 2.790736
$ws.Range("N4").Value = 8.372208000000001
$ws.Range("O4").Value = 0.0329987043561157
$ws.Range("P4").Value = 0.0329987043561157
$ws.Range("Q4").Value = 33.39717957853333
$ws.Range("R4").Value = 300.5746162068
$ws.Range("S4").Value = 0.02710311965497586
$ws.Range("T4").Value = 0.02710311965497585
$ws.Range("G5").Value = 11.96715833333333
$ws.Range("H5").Value = 35.901475
$ws.Range("I5").Value = 0.8213389035667636
$ws.Range("J5").Value = 0.8213389035667636
$ws.Range("M5").Value = 15.50544933333333
$ws.Range("N5").Value = 46.516348
$ws.Range("O5").Value = 0.1833422217147727
$ws.Range("P5").Value = 0.1833422217147727
$ws.Range("Q5").Value = 185.5561672014778
$ws.Range("R5").Value = 1670.0055048133
$ws.Range("S5").Value = 0.1505860993607059
$ws.Range("T5").Value = 0.1505860993607059
$ws.Range("I6").Value = 0.008600300405516565
$ws.Range("J6").Value = 0.008600300405516565
$ws.Range("M6").Value = 57.65261933333333
$ws.Range("N6").Value = 172.957858
$ws.Range("O6").Value = 0.6817060950001529
$ws.Range("P6").Value = 0.6817060950001529
$ws.Range("Q6").Value = 7.224392076040666
$ws.Range("R6").Value = 65.01952868436599
$ws.Range("S6").Value = 0.005862877205272929
$ws.Range("T6").Value = 0.005862877205272929
$ws.Range("I7").Value = 0.008600300405516565
$ws.Range("J7").Value = 0.008600300405516565
$ws.Range("O7").Value = 0.1019529789289588
$ws.Range("P7").Value = 0.1019529789289588
$ws.Range("S7").Value = 0.0008768262460263464
$ws.Range("T7").Value = 0.0008768262460263464
$ws.Range("I8").Value = 0.008600300405516565
$ws.Range("J8").Value = 0.008600300405516565
$ws.Range("M8").Value = 2.790736
$ws.Range("N8").Value = 8.372208000000001
$ws.Range("O8").Value = 0.0329987043561157
$ws.Range("P8").Value = 0.0329987043561157
$ws.Range("Q8").Value = 0.349704337424
$ws.Range("R8").Value = 3.147339036816
$ws.Range("S8").Value = 0.0002837987704554231
$ws.Range("T8").Value = 0.0002837987704554231
$ws.Range("I9").Value = 0.008600300405516565
$ws.Range("J9").Value = 0.008600300405516565
$ws.Range("M9").Value = 15.50544933333333
$ws.Range("N9").Value = 46.516348
$ws.Range("O9").Value = 0.1833422217147727
$ws.Range("P9").Value = 0.1833422217147727
$ws.Range("Q9").Value = 1.942972350510667
$ws.Range("R9").Value = 17.486751154596
$ws.Range("S9").Value = 0.001576798183761867
$ws.Range("T9").Value = 0.001576798183761867
$ws.Range("G10").Value = 2.477837666666666
$ws.Range("H10").Value = 7.433513
$ws.Range("I10").Value = 0.1700607960277199
$ws.Range("J10").Value = 0.1700607960277199
$ws.Range("M10").Value = 57.65261933333333
$ws.Range("N10").Value = 172.957858
$ws.Range("O10").Value = 0.6817060950001529
$ws.Range("P10").Value = 0.6817060950001529
$ws.Range("Q10").Value = 142.8538317661282
$ws.Range("R10").Value = 1285.684485895154
$ws.Range("S10").Value = 0.1159314811726744
$ws.Range("T10").Value = 0.1159314811726744
$ws.Range("G11").Value = 2.477837666666666
$ws.Range("H11").Value = 7.433513
$ws.Range("I11").Value = 0.1700607960277199
$ws.Range("J11").Value = 0.1700607960277199
$ws.Range("O11").Value = 0.1019529789289588
$ws.Range("P11").Value = 0.1019529789289588
$ws.Range("Q11").Value = 21.36459363762888
$ws.Range("R11").Value = 192.28134273866
$ws.Range("S11").Value = 0.01733820475405609
$ws.Range("T11").Value = 0.01733820475405609
$ws.Range("G12").Value = 2.477837666666666
$ws.Range("H12").Value = 7.433513
$ws.Range("I12").Value = 0.1700607960277199
$ws.Range("J12").Value = 0.1700607960277199
$ws.Range("M12").Value = 2.790736
$ws.Range("N12").Value = 8.372208000000001
$ws.Range("O12").Value = 0.0329987043561157
$ws.Range("P12").Value = 0.0329987043561157
$ws.Range("Q12").Value = 6.914990778522666
$ws.Range("R12").Value = 62.234917006704
$ws.Range("S12").Value = 0.005611785930684424
$ws.Range("T12").Value = 0.005611785930684422
$ws.Range("G13").Value = 2.477837666666666
$ws.Range("H13").Value = 7.433513
$ws.Range("I13").Value = 0.1700607960277199
$ws.Range("J13").Value = 0.1700607960277199
$ws.Range("M13").Value = 15.50544933333333
$ws.Range("N13").Value = 46.516348
$ws.Range("O13").Value = 0.1833422217147727
$ws.Range("P13").Value = 0.1833422217147727
$ws.Range("Q13").Value = 38.41998639672488
$ws.Range("R13").Value = 345.779877570524
$ws.Range("S13").Value = 0.03117932417030495
$ws.Range("T13").Value = 0.03117932417030495
